$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename sheet to match the re-exported MySQL-for-Excel workbook.
$ws.Name = "Sheet1"

# Register the hidden MySQL-for-Excel helper defined name that the
# add-in writes into every workbook it touches.
$n = $wb.Names.Add("LOCAL_MYSQL_DATE_FORMAT", "=REPT(LOCAL_YEAR_FORMAT,4)&LOCAL_DATE_SEPARATOR&REPT(LOCAL_MONTH_FORMAT,2)&LOCAL_DATE_SEPARATOR&REPT(LOCAL_DAY_FORMAT,2)&"" ""&REPT(LOCAL_HOUR_FORMAT,2)&LOCAL_TIME_SEPARATOR&REPT(LOCAL_MINUTE_FORMAT,2)&LOCAL_TIME_SEPARATOR&REPT(LOCAL_SECOND_FORMAT,2)")
$n.Visible = $false
